# Update the "K" column (column G) values in the save-data sheet.
# The commit regenerates the save_data to use K instead of Strike#, recomputing
# std/mean and writing the new s_vals into column G for rows 2-27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newKValues = @{
    2  = 1
    3  = 9
    4  = 8
    5  = 5
    6  = 5
    7  = 4
    8  = 7
    9  = 5
    10 = 4
    11 = 5
    12 = 7
    13 = 2
    14 = 5
    15 = 2
    16 = 3
    17 = 5
    18 = 8
    19 = 3
    20 = 4
    21 = 4
    22 = 4
    23 = 3
    24 = 3
    25 = 3
    26 = 6
    27 = 1
}

foreach ($row in $newKValues.Keys) {
    $ws.Range("G$row").Value = $newKValues[$row]
}
